$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells to lowercase language-code variants
$ws.Range("B1").Value = "Description (en)"
$ws.Range("C1").Value = "Description (pl)"

# Update the selected/active cell shown when the sheet is reopened
$ws.Range("B2").Select()
